$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOM3241")

# Name: Chemistry of Materials -> Materials chemistry
$ws.Range("B4").Value = "Materials chemistry"
$ws.Range("C4").Value = "Materials chemistry"

# Ativacao date: 01/01/2012 -> 01/01/2023 (both occurrences)
$ws.Range("B8").Value = "01/01/2023"
$ws.Range("C8").Value = "01/01/2023"
$ws.Range("B13").Value = "01/01/2023"
$ws.Range("C13").Value = "01/01/2023"

# New Objectives text (row 11, next to "Objectives:")
$ws.Range("B11").Value = "Provide the student with the main types of organic and inorganic synthesis of materials as well as presenting the main analytical techniques for material characterization."
$ws.Range("C11").Value = "Provide the student with the main types of organic and inorganic synthesis of materials as well as presenting the main analytical techniques for material characterization."
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Style = "Normal"

# New Short syllabus text (row 14, next to "Short syllabus:")
$ws.Range("B14").Value = "Introduction to the chemistry of materials and its association with the synthesis of new materials. The modern view of the atom and chemical bonds. Crystal structure and crystallographic characterization techniques. Epitaxial thin films and films in general and their impact on modern technology. Amorphous materials, synthesis and applications. Synthesis of materials and chemical transformations. Processes and Techniques of crystal growth in general. Conducting polymers and their applications in modern technology."
$ws.Range("C14").Value = "Introduction to the chemistry of materials and its association with the synthesis of new materials. The modern view of the atom and chemical bonds. Crystal structure and crystallographic characterization techniques. Epitaxial thin films and films in general and their impact on modern technology. Amorphous materials, synthesis and applications. Synthesis of materials and chemical transformations. Processes and Techniques of crystal growth in general. Conducting polymers and their applications in modern technology."

# New Syllabus text (row 16, next to "Syllabus:")
$ws.Range("B16").Value = "Materials chemistry: definition; role of chemistry in materials science; fundamentals.Atomistics and the modern view of the atom with quantum foundations.Types of chemical bonds: van der Waals forces, Lennard-Jones potential, covalent bonding, coordination bonds, ionic bonds and metallic bonds.Polycrystalline and monocrystalline materials. The crystallographic order and crystallographic and microscopic characterization techniques. The importance of single crystals in electronic applications. High quality crystal growth techniques such as: flow method, Czochralski method, Brigdmann method, vapor transport method and modified isothermal vapor transport growth method. Amorphous materials and their importance for modern technology. Concepts and techniques for growing amorphous materials. Epitaxial thin films, growth techniques such as: chemical vapor, sputtering, laser ablation and MBE. Thin films grown by electrolysis for protective coating, concepts and applications. Synthesis of conductive polymers, concepts and applications as electronic devices."
$ws.Range("C16").Value = "Materials chemistry: definition; role of chemistry in materials science; fundamentals.Atomistics and the modern view of the atom with quantum foundations.Types of chemical bonds: van der Waals forces, Lennard-Jones potential, covalent bonding, coordination bonds, ionic bonds and metallic bonds.Polycrystalline and monocrystalline materials. The crystallographic order and crystallographic and microscopic characterization techniques. The importance of single crystals in electronic applications. High quality crystal growth techniques such as: flow method, Czochralski method, Brigdmann method, vapor transport method and modified isothermal vapor transport growth method. Amorphous materials and their importance for modern technology. Concepts and techniques for growing amorphous materials. Epitaxial thin films, growth techniques such as: chemical vapor, sputtering, laser ablation and MBE. Thin films grown by electrolysis for protective coating, concepts and applications. Synthesis of conductive polymers, concepts and applications as electronic devices."

# Norma de recuperacao text change
$ws.Range("B20").Value = "Média simples de duas provas escritas,  Conceito Final = (P1 + P2)/2"
$ws.Range("C20").Value = "Média simples de duas provas escritas,  Conceito Final = (P1 + P2)/2"

# Bibliografia text change
$ws.Range("B21").Value = "Aplicação de duas provas escritas dentro do prazo regimental antes do início do próximo semestre letivo."
$ws.Range("C21").Value = "Aplicação de duas provas escritas dentro do prazo regimental antes do início do próximo semestre letivo."
